# DAX and EP Global Objects
# - Remove the unused "Range" Param rows (fromRow/fromCol/toRow/toCol) on the
#   RVL sheet (old rows 17-20), which shifts everything below up by 4 rows.
# - Rename the "Functions" object / "Ep*" actions to the new "EP" object with
#   the "Ep" prefix stripped from the action name (Launch, ChangeCompany,
#   OpenModule, SelectTab, ClickRibbon, FilterGrid, SelectFastTab) on both
#   the RVL and Cleanup sheets.

$wb = $excel.ActiveWorkbook

$rvl = $wb.Worksheets.Item("RVL")

# Remove the 4 now-unused Param rows (fromRow/fromCol/toRow/toCol) under the
# "Range" map block.
$rvl.Range("A17:A20").EntireRow.Delete() | Out-Null

# Walk every row and fix up any "Functions"/"Ep*" pair into "EP"/"*".
$usedRows = $rvl.UsedRange.Rows.Count
for ($r = 1; $r -le $usedRows; $r++) {
    $objectCell = $rvl.Cells.Item($r, 3)
    $actionCell = $rvl.Cells.Item($r, 4)
    $objectVal = $objectCell.Value2
    $actionVal = $actionCell.Value2
    if ($objectVal -eq "Functions" -and $actionVal -ne $null -and $actionVal.StartsWith("Ep")) {
        $objectCell.Value2 = "EP"
        $actionCell.Value2 = $actionVal.Substring(2)
    }
}

$cleanup = $wb.Worksheets.Item("Cleanup")
$usedRowsCleanup = $cleanup.UsedRange.Rows.Count
for ($r = 1; $r -le $usedRowsCleanup; $r++) {
    $objectCell = $cleanup.Cells.Item($r, 3)
    $actionCell = $cleanup.Cells.Item($r, 4)
    $objectVal = $objectCell.Value2
    $actionVal = $actionCell.Value2
    if ($objectVal -eq "Functions" -and $actionVal -ne $null -and $actionVal.StartsWith("Ep")) {
        $objectCell.Value2 = "EP"
        $actionCell.Value2 = $actionVal.Substring(2)
    }
}
